$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Row 3 - existing "Wages Expense" (B3) entry also gets a Debit amount now
$ws.Range("D3").Formula = "=60000+260000"

# Row 4 - new entry: TAX - IURAN ARIESTA
$ws.Range("B4").Value = "TAX - IURAN ARIESTA"
$ws.Range("D4").Value = 660000

# Row 5 - new entry: TAX - P.Tata
$ws.Range("B5").Value = "TAX - P.Tata"
$ws.Range("D5").Value = 200000

# Row 6 - new entry: A/R
$ws.Range("B6").Value = "A/R"
$ws.Range("C6").Formula = "=6420000+5580000+5000000+1794000+1900000+70000000+29598000"

# Row 7 - new entry: TRANSFER BCA
$ws.Range("B7").Value = "TRANSFER BCA"
$ws.Range("D7").Formula = "=1787000+1794000"

# Row 8 - new entry: SALES - cash/retail
$ws.Range("B8").Value = "SALES - cash/retail"
$ws.Range("C8").Formula = "=124499525-87127525-29598000"

# Row 9 - new entry: SELISIH - lebih
$ws.Range("B9").Value = "SELISIH - lebih"
$ws.Range("C9").Value = 446500

# Row 10 - new entry: SETOR KE BANK
$ws.Range("B10").Value = "SETOR KE BANK"
$ws.Range("D10").Value = 124000000

# Row 11 - new day (9-Feb-2021), Wages Expense entry
$ws.Range("A11").Value = 44236
$ws.Range("B11").Value = "Wages Expense"
$ws.Range("D11").Formula = "=60000"

# Row 12 - TRANSFER BCA
$ws.Range("B12").Value = "TRANSFER BCA"
$ws.Range("D12").Formula = "=8839000+6709000+3720000+842000+25000000"

# Row 13 - A/R
$ws.Range("B13").Value = "A/R"
$ws.Range("C13").Formula = "=6000000"

# Update the active selection to C7, matching the author's last cursor position
$ws.Range("C7").Select()
